# Completed code and plotting for VT 1.
#
# - "General" sheet: splits the old "Solar Degradation (%/year)" input row
#   into two separate rows - "ST Degradation (%/year)" and
#   "PV Degradation (%/year)" - pushing "Discount Rate (%)" down a row.
# - "Fuel" sheet becomes the active tab (was "ST Pricing").

$wb = $excel.ActiveWorkbook

$wsGeneral = $wb.Worksheets.Item("General")
$wsFuel = $wb.Worksheets.Item("Fuel")
$wsST = $wb.Worksheets.Item("ST Pricing")

# Insert a new row 4 on "General" (below the old "Solar Degradation" row)
# for the new "ST Degradation (%/year)" input, then relabel row 3 as the
# "PV Degradation (%/year)" input. This ordering keeps the shared-string
# table append order (ST Degradation before PV Degradation) matching the
# authored workbook.
$wsGeneral.Rows.Item(4).Insert() | Out-Null
$wsGeneral.Range("A4").Value = "ST Degradation (%/year)"
$wsGeneral.Range("B4").Value = 2
$wsGeneral.Range("A3").Value = "PV Degradation (%/year)"

# Move the selection on "General" down to the now-relocated "Discount Rate" row
$wsGeneral.Range("A5").Select() | Out-Null

# "ST Pricing" keeps its own cursor position even though it's no longer the
# active tab.
$wsST.Range("C14").Select() | Out-Null

# "Fuel" becomes the active sheet/tab.
$wsFuel.Activate() | Out-Null
